$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.560.24"
$ws.Range("E2").Value = "  -1.20%  "
# Row 3
$ws.Range("D3").Value = "1.593.83"
$ws.Range("E3").Value = "  -2.24%  "
# Row 4
$ws.Range("E4").Value = "  +0.37%  "
# Row 5
$ws.Range("D5").Value = "'207.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.80%  "
# Row 6
$ws.Range("E6").Value = "  -3.35%  "
# Row 7
$ws.Range("E7").Value = "  +0.45%  "
# Row 8
$ws.Range("D8").Value = "'22.29"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.89%  "
# Row 9
$ws.Range("E9").Value = "  -2.07%  "
# Row 10
$ws.Range("E10").Value = "  -3.41%  "
# Row 11
$ws.Range("D11").Value = "'0.0865"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.87%  "
# Row 12
$ws.Range("D12").Value = "1.822.38"
$ws.Range("E12").Value = "  -2.12%  "
# Row 13
$ws.Range("D13").Value = "1.605.62"
$ws.Range("E13").Value = "  -1.54%  "
# Row 14
$ws.Range("D14").Value = "'3.86"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.03%  "
# Row 15
$ws.Range("D15").Value = "'0.536"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.47%  "
# Row 16
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'63.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.17%  "
# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "27.559.77"
$ws.Range("E17").Value = "  -1.19%  "
# Row 18
$ws.Range("D18").Value = "'217.26"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.25%  "
# Row 19
$ws.Range("D19").Value = "'7.37"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.98%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0692"
$ws.Range("E20").Value = "  -3.82%  "
# Row 21
$ws.Range("E21").Value = "  +0.36%  "
# Row 22
$ws.Range("D22").Value = "'4.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.69%  "
# Row 23
$ws.Range("D23").Value = "'9.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.78%  "
# Row 24
$ws.Range("E24").Value = "  -3.45%  "
# Row 25
$ws.Range("D25").Value = "'152.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.36%  "
# Row 26
$ws.Range("D26").Value = "'6.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.74%  "
# Row 27
$ws.Range("E27").Value = "  +0.40%  "
# Row 29
$ws.Range("E29").Value = "  -3.96%  "
# Row 30
$ws.Range("E30").Value = "  -1.85%  "
# Row 31
$ws.Range("E31").Value = "  -3.08%  "
# Row 32
$ws.Range("D32").Value = "'3.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.33%  "
# Row 33
$ws.Range("D33").Value = "1.373.86"
$ws.Range("E33").Value = "  -1.13%  "
# Row 34
$ws.Range("E34").Value = "  -5.34%  "
# Row 35
$ws.Range("E35").Value = "  -3.97%  "
# Row 36
$ws.Range("E36").Value = "  -5.52%  "
# Row 37
$ws.Range("E37").Value = "  -1.26%  "
# Row 38
$ws.Range("E38").Value = "  -3.38%  "
# Row 39
$ws.Range("D39").Value = "'0.539"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.09%  "
# Row 40
$ws.Range("E40").Value = "  -4.32%  "
# Row 41
$ws.Range("E41").Value = "  +0.39%  "
# Row 43
$ws.Range("E43").Value = "  +2.86%  "
# Row 44
$ws.Range("D44").Value = "'5.35"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.27%  "
# Row 45
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'63.96"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.49%  "
# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.07%  "
# Row 47
$ws.Range("D47").Value = "1.732.66"
$ws.Range("E47").Value = "  -2.23%  "
# Row 48
$ws.Range("D48").Value = "'87.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.59%  "
# Row 49
$ws.Range("E49").Value = "  -3.07%  "
# Row 50
$ws.Range("D50").Value = "'0.0970"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.39%  "
# Row 51
$ws.Range("D51").Value = "'0.0497"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
